# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Vega Modelo de Temuco - Pepino dulce"
# at row 321, pushing the existing rows 321:335 down to 322:336.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 321:335 down one row, creating a blank row 321.
$ws.Rows.Item(321).Insert()

# Populate the newly inserted row 321 with the new record's data.
$ws.Range("A321").Value = 10
$ws.Range("B321").Value = "Vega Modelo de Temuco"
$ws.Range("C321").Value = "La Araucanía"
$ws.Range("D321").Value = 45041
$ws.Range("E321").Value = 9
$ws.Range("F321").Value = 100112043
$ws.Range("G321").Value = "Pepino dulce"
$ws.Range("H321").Value = "Cultivar IV Región"
$ws.Range("I321").Value = "Primera"
$ws.Range("J321").Value = 260
$ws.Range("K321").Value = 17000
$ws.Range("L321").Value = 18000
$ws.Range("M321").Value = 17462
$ws.Range("N321").Value = "$/bandeja 18 kilos"
$ws.Range("O321").Value = "Provincia de Limarí"
$ws.Range("P321").Value = 970
$ws.Range("Q321").Value = 18
$ws.Range("R321").Value = "Hortaliza"
